$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fill in the real spell data for FERDINAND / Dart (was GUNNAR / Test spell placeholder) ---
$ws.Range("A2").Value = "FERDINAND"
$ws.Range("B2").Value = "Dart"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 43
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "WFRP, Page 240"
$ws.Range("H2").Value = "No description"

# --- Row 3: new spell, FERDINAND / Light ---
$ws.Range("A3").Value = "FERDINAND"
$ws.Range("B3").Value = "Light"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -1
$ws.Range("E3").Value = -1
$ws.Range("F3").Value = 2580
$ws.Range("G3").Value = "WFRP, Page 241"
$ws.Range("G3").Font.Bold = $true
$ws.Range("H3").Value = "No description"

# --- Row 4: new spell, FERDINAND / Shock ---
$ws.Range("A4").Value = "FERDINAND"
$ws.Range("B4").Value = "Shock"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = "WFRP, Page 242"
$ws.Range("G4").Font.Bold = $true
$ws.Range("H4").Value = "No description"

# --- Column widths: widen column A, and size the two new columns G/H ---
$ws.Columns.Item(1).ColumnWidth = 12.25
$ws.Columns.Item(7).ColumnWidth = 18.25
$ws.Columns.Item(8).ColumnWidth = 13.75

# --- Selection moves to A8 ---
[void]$ws.Range("A8").Select()
